$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 4 new rows for D_TYPE ::= Bool_Ni / Double_Nil / Int_Nil / String_Nil
#        right after the existing D_TYPE ::= String row (row 17), i.e. before row 18.
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = "D_TYPE"
$ws.Range("B18").Formula = '=$I$2'
$ws.Range("C18").Value = "Bool_Ni"

$ws.Range("A19").Value = "D_TYPE"
$ws.Range("B19").Formula = '=$I$2'
$ws.Range("C19").Value = "Double_Nil"

$ws.Range("A20").Value = "D_TYPE"
$ws.Range("B20").Formula = '=$I$2'
$ws.Range("C20").Value = "Int_Nil"

$ws.Range("A21").Value = "D_TYPE"
$ws.Range("B21").Formula = '=$I$2'
$ws.Range("C21").Value = "String_Nil"

# --- 2. Update the two "if EXP { ... }" productions to use the new IF_COND non-terminal.
#        After the insert above, these productions now live at rows 41 and 46.
$ws.Range("C41").Value = "if IF_COND { FUNC_STMT_LIST } FUNC_ELSE_CLAUSE"
$ws.Range("C46").Value = "if IF_COND { STMT_LIST } ELSE_CLAUSE"

# --- 3. Append the two new IF_COND productions at the end of the grammar table.
$ws.Range("A63").Value = "IF_COND"
$ws.Range("B63").Formula = '=$I$2'
$ws.Range("C63").Value = "EXP"

$ws.Range("A64").Value = "IF_COND"
$ws.Range("B64").Formula = '=$I$2'
$ws.Range("C64").Value = "let id"

# --- 4. Restore the active selection / view the author left the sheet in.
$ws.Range("C21").Select()
